$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" banner timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 01:50"

# Re-sequence country labels in column A (rows whose country moved position
# in the published list order) to match the updated source data ordering.
$ws.Range("A147").Value = "Mali"   # was "Niger"
$ws.Range("A148").Value = "Niger"   # was "Mali"
$ws.Range("A155").Value = "Bahamas"   # was "Mongolia"
$ws.Range("A156").Value = "Eritrea"   # was "Guinea Ecuatorial"
$ws.Range("A157").Value = "Guinea Ecuatorial"   # was "Eritrea"
$ws.Range("A158").Value = "Mongolia"   # was "San Martin (Parte Francesa)"
$ws.Range("A159").Value = "Islas Caimanes"   # was "Dominica"
$ws.Range("A160").Value = "Dominica"   # was "Bahamas"
$ws.Range("A161").Value = "San Martin (Parte Francesa)"   # was "Namibia"
$ws.Range("A162").Value = "Namibia"   # was "Birmania"
$ws.Range("A163").Value = "Birmania"   # was "Groenlandia"
$ws.Range("A164").Value = "Groenlandia"   # was "Suazilandia"
$ws.Range("A165").Value = "Suazilandia"   # was "Granada"
$ws.Range("A166").Value = "Granada"   # was "Siria"
$ws.Range("A168").Value = "Siria"   # was "Laos"
$ws.Range("A169").Value = "Laos"   # was "Seychelles"
$ws.Range("A170").Value = "Seychelles"   # was "Surinam"
$ws.Range("A171").Value = "Surinam"   # was "Mozambique"
$ws.Range("A172").Value = "Mozambique"   # was "Libia"
$ws.Range("A173").Value = "Libia"   # was "Guyana"
$ws.Range("A174").Value = "Guyana"   # was "Islas Caimanes"
$ws.Range("A177").Value = "Zimbabue"   # was "Gabon"
$ws.Range("A178").Value = "Gabon"   # was "Zimbabue"
$ws.Range("A180").Value = "Santa Sede"   # was "Benin"
$ws.Range("A182").Value = "Benin"   # was "Santa Sede"
$ws.Range("A183").Value = "Sudan"   # was "Cabo Verde"
$ws.Range("A184").Value = "Cabo Verde"   # was "Sudan"
$ws.Range("A192").Value = "Nicaragua"   # was "Gambia"
$ws.Range("A193").Value = "Gambia"   # was "Nicaragua"

# Apply updated numeric statistics (Casos totales / Nuevos casos / Casos
# activos / Recuperados / Casos criticos / Muertes hoy / Muertes).
$ws.Range("B4").Value = 141781
$ws.Range("C4").Value = 18203
$ws.Range("E4").Value = 134875
$ws.Range("G4").Value = 251
$ws.Range("H4").Value = 2471
$ws.Range("D18").Value = 573
$ws.Range("E18").Value = 5642
$ws.Range("E43").Value = 1065
$ws.Range("G43").Value = 7
$ws.Range("H43").Value = 39
$ws.Range("F46").Value = 36
$ws.Range("E80").Value = 256
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 3
$ws.Range("E88").Value = 193
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 6
$ws.Range("F93").Value = 3
$ws.Range("C147").Value = 0
$ws.Range("C148").Value = 8
$ws.Range("B155").Value = 14
$ws.Range("C155").Value = 4
$ws.Range("D155").Value = 1
$ws.Range("E155").Value = 13
$ws.Range("C156").Value = 6
$ws.Range("C157").Value = 0
$ws.Range("B158").Value = 12
$ws.Range("E158").Value = 12
$ws.Range("B159").Value = 12
$ws.Range("C159").Value = 4
$ws.Range("H159").Value = 1
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 11
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 11
$ws.Range("B162").Value = 11
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 2
$ws.Range("E162").Value = 9
$ws.Range("C163").Value = 2
$ws.Range("D163").Value = 0
$ws.Range("E163").Value = 10
$ws.Range("B164").Value = 10
$ws.Range("D164").Value = 2
$ws.Range("E164").Value = 8
$ws.Range("C165").Value = 0
$ws.Range("C166").Value = 2
$ws.Range("E166").Value = 9
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
$ws.Range("B168").Value = 9
$ws.Range("C168").Value = 4
$ws.Range("G168").Value = 1
$ws.Range("H168").Value = 1
$ws.Range("C172").Value = 0
$ws.Range("C173").Value = 5
$ws.Range("E173").Value = 8
$ws.Range("H173").Value = 0
$ws.Range("C183").Value = 1
$ws.Range("C184").Value = 0
$ws.Range("C192").Value = 0
$ws.Range("C193").Value = 1
